# Apply the "Hillclimber" / "Simulated annealing" worksheet split described
# by the target diff.

$wb = $excel.ActiveWorkbook

# --- 1. Existing sheet: add the two new hillclimber rows (11 & 12), then
#        rename it "Hillclimber" -----------------------------------------
$hill = $wb.Worksheets.Item(1)

$hill.Cells.Item(11, 1).Value = 500
$hill.Cells.Item(11, 2).Value = 800
$hill.Cells.Item(11, 3).Value = 80

$hill.Cells.Item(12, 1).Value = 1000
$hill.Cells.Item(12, 2).Value = 800
$hill.Cells.Item(12, 3).Value = 80

$hill.Name = "Hillclimber"

# --- 2. New worksheet, inserted right after "Hillclimber" ----------------
$sa = $wb.Worksheets.Add($null, $hill)
$sa.Name = "Simulated annealing"

# Header row. Shared-string table order matters for an exact match: the
# source file has "...temperature" (index 32) created before
# "...repeats" (index 33), so write B1 before A1.
$sa.Cells.Item(1, 2).Value  = "simulated annealing lessons -> temperature"
$sa.Cells.Item(1, 1).Value  = "simulated annealing lessons -> repeats"
$sa.Cells.Item(1, 3).Value  = "hillclimber students outer -> 0"
$sa.Cells.Item(1, 4).Value  = "hillclimber students inner -> i"
$sa.Cells.Item(1, 5).Value  = "runtime"
$sa.Cells.Item(1, 6).Value  = "average"
$sa.Cells.Item(1, 7).Value  = "minimum"
$sa.Cells.Item(1, 8).Value  = "intermediate minimum"
$sa.Cells.Item(1, 9).Value  = "pickle"
$sa.Cells.Item(1, 10).Value = "schedule"

# Data rows. Column B ("temperature") holds text that looks numeric
# ("0.5" / "1.0" / "2.0"), so force a text number format before writing it
# and then clear the formatting again so the written file doesn't carry an
# explicit style index on those cells (matches the source workbook, which
# was produced without any extra per-cell formatting).
$tempRange = $sa.Range("B2:B7")
$tempRange.NumberFormat = "@"

$sa.Cells.Item(2, 1).Value = 10000
$sa.Cells.Item(2, 2).Value = "0.5"
$sa.Cells.Item(2, 3).Value = 300
$sa.Cells.Item(2, 4).Value = 50

$sa.Cells.Item(3, 1).Value = 10000
$sa.Cells.Item(3, 2).Value = "1.0"
$sa.Cells.Item(3, 3).Value = 300
$sa.Cells.Item(3, 4).Value = 50

$sa.Cells.Item(4, 1).Value = 10000
$sa.Cells.Item(4, 2).Value = "2.0"
$sa.Cells.Item(4, 3).Value = 300
$sa.Cells.Item(4, 4).Value = 50

$sa.Cells.Item(5, 1).Value = 30000
$sa.Cells.Item(5, 2).Value = "0.5"
$sa.Cells.Item(5, 3).Value = 300
$sa.Cells.Item(5, 4).Value = 50

$sa.Cells.Item(6, 1).Value = 30000
$sa.Cells.Item(6, 2).Value = "1.0"
$sa.Cells.Item(6, 3).Value = 300
$sa.Cells.Item(6, 4).Value = 50

$sa.Cells.Item(7, 1).Value = 30000
$sa.Cells.Item(7, 2).Value = "2.0"
$sa.Cells.Item(7, 3).Value = 300
$sa.Cells.Item(7, 4).Value = 50

$tempRange.ClearFormats()

# Column widths (best-effort match of the authored "best fit" widths).
$sa.Columns.Item(1).ColumnWidth  = 31.8307291666667
$sa.Columns.Item(2).ColumnWidth  = 36.0533854166667
$sa.Columns.Item(3).ColumnWidth  = 24.7213541666667
$sa.Columns.Item(4).ColumnWidth  = 23.7213541666667
$sa.Columns.Item(5).ColumnWidth  = 6.49869791666667
$sa.Columns.Item(6).ColumnWidth  = 6.60807291666667
$sa.Columns.Item(7).ColumnWidth  = 7.72135416666667
$sa.Columns.Item(8).ColumnWidth  = 18.8307291666667
$sa.Columns.Item(9).ColumnWidth  = 4.83072916666667
$sa.Columns.Item(10).ColumnWidth = 7.27604166666667

# --- 3. Selections -------------------------------------------------------
# "Hillclimber" keeps a (non-active) selection of A16:I16 ...
$hill.Range("A16:I16").Select()

# ... while "Simulated annealing" ends up the active tab, selection C7:D7.
$sa.Activate()
$sa.Range("C7:D7").Select()
